# Updates the LR-pairs sheet with new TPM-derived values, adding the "ECs"
# sending/target cluster rows and refreshing all computed statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full table (header row already present in A1:T1). Data rows 2-10 cover all
# combinations of Sending cluster x Target cluster in {ECs, FAPs, MuSCs}
# for Ligand=Angpt4 / Receptor=Tek.
$data = @(
    @("ECs",   "Angpt4", "Tek", "ECs",   2, 0.6666666666666666, 0.04873833333333333, 0.146215,   0.02106797933705496, 0.02106797933705496, 3, 1, 43.28265566666667, 129.847967, 0.667219228070094,  0.667219228070094,  2.109524499433889,  18.985720494905,    0.0140569609102665,  0.0140569609102665),
    @("ECs",   "Angpt4", "Tek", "FAPs",  2, 0.6666666666666666, 0.04873833333333333, 0.146215,   0.02106797933705496, 0.02106797933705496, 3, 1, 18.63243533333333, 55.897306,  0.2872263480299067, 0.2872263480299067, 0.9081138440877777, 8.173024596789999,  0.006051278765351829, 0.006051278765351829),
    @("ECs",   "Angpt4", "Tek", "MuSCs", 2, 0.6666666666666666, 0.04873833333333333, 0.146215,   0.02106797933705496, 0.02106797933705496, 3, 1, 2.955125333333334, 8.865376000000001, 0.04555442389999943, 0.04555442389999944, 0.1440278835377778, 1.29625095184,      0.0009597396614366305, 0.0009597396614366307),
    @("FAPs",  "Angpt4", "Tek", "ECs",   3, 1,                  1.741645,            5.224935,   0.7528558808429727,  0.7528558808429727,  3, 1, 43.28265566666667, 129.847967, 0.667219228070094,  0.667219228070094,  75.38302082857167,  678.4471874571451,  0.502319919664079,   0.502319919664079),
    @("FAPs",  "Angpt4", "Tek", "FAPs",  3, 1,                  1.741645,            5.224935,   0.7528558808429727,  0.7528558808429727,  3, 1, 18.63243533333333, 55.897306,  0.2872263480299067, 0.2872263480299067, 32.45108783612334,  292.05979052511,    0.2162400452473656,  0.2162400452473656),
    @("FAPs",  "Angpt4", "Tek", "MuSCs", 3, 1,                  1.741645,            5.224935,   0.7528558808429727,  0.7528558808429727,  3, 1, 2.955125333333334, 8.865376000000001, 0.04555442389999943, 0.04555442389999944, 5.146779261173334,  46.32101335056001,  0.03429591593152824, 0.03429591593152825),
    @("MuSCs", "Angpt4", "Tek", "ECs",   3, 1,                  0.5230009999999999, 1.569003,   0.2260761398199722,  0.2260761398199722,  3, 1, 43.28265566666667, 129.847967, 0.667219228070094,  0.667219228070094,  22.63687219632233,  203.731849766901,   0.1508423474957485, 0.1508423474957485),
    @("MuSCs", "Angpt4", "Tek", "FAPs",  3, 1,                  0.5230009999999999, 1.569003,   0.2260761398199722,  0.2260761398199722,  3, 1, 18.63243533333333, 55.897306,  0.2872263480299067, 0.2872263480299067, 9.744782311768665,  87.70304080591799,  0.06493502401718919, 0.06493502401718919),
    @("MuSCs", "Angpt4", "Tek", "MuSCs", 3, 1,                  0.5230009999999999, 1.569003,   0.2260761398199722,  0.2260761398199722,  3, 1, 2.955125333333334, 8.865376000000001, 0.04555442389999943, 0.04555442389999944, 1.545533504458667,  13.909801540128,    0.01029876830703456, 0.01029876830703456)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($value in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $value
        $colIndex++
    }
    $rowIndex++
}
